$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the " " run and the "A cluster of heavy industrial drums, " run
#    (in the final narrative "hiding spots" paragraph, not the numbered-list
#    one earlier in the document) into a single run.
# ---------------------------------------------------------------------------
$searchRange = $d.Content
$found = $searchRange.Find.Execute("space.\n\n A cluster of heavy industrial drums, ", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $prefixLen = ("space.\n\n").Length
    $targetStart = $searchRange.Start + $prefixLen
    $targetEnd = $searchRange.End
    $mergeRange = $d.Range($targetStart, $targetEnd)
    $mergeRange.Find.Execute(" A cluster of heavy industrial drums, ", $false, $false, $false, `
        $false, $false, $true, 1, $false, " A cluster of heavy industrial drums, ", 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 2) After "...between them and the wall." insert a new blank paragraph and a
#    new narrative paragraph describing the arachnid creature.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("between them and the wall.", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$anchor.Collapse(0)
$insertAt = $d.Range($anchor.End, $anchor.End)

$insertAt.Find.Execute("between them and the wall.", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null

# Locate the *last* occurrence of the paragraph-ending phrase (the narrative
# "hiding spots" paragraph near the end of the document), not the numbered
# list entry earlier in the doc.
$searchCursor = $d.Content
$lastStart = -1
$lastEnd = -1
while ($searchCursor.Find.Execute("between them and the wall.", $false, $false, $false, `
        $false, $false, $true, 1, $false, "", 0)) {
    $lastStart = $searchCursor.Start
    $lastEnd = $searchCursor.End
    $searchCursor.Collapse(0)
    $searchCursor.End = $d.Content.End
}

$insertionPoint = $d.Range($lastEnd, $lastEnd)

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$blankParaXml = "<w:p $wNs><w:pPr><w:jc w:val='left'/><w:rPr><w:lang w:val='en-FI'/></w:rPr></w:pPr></w:p>"

$newParaXml = @"
<w:p $wNs>
  <w:pPr>
    <w:jc w:val="left"/>
    <w:rPr><w:lang w:val="en-FI"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:lang w:val="en-FI"/></w:rPr>
    <w:t xml:space="preserve">    The source of the clicking is revealed: a colossal arachnid creature drops from a vent high above the console array and lands heavily on the floor, its spike-like legs scattering </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr><w:lang w:val="en-FI"/></w:rPr>
    <w:t>debris.\n\n</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr><w:lang w:val="en-FI"/></w:rPr>
    <w:t>&quot;</w:t>
  </w:r>
</w:p>
"@

$insertionPoint.InsertXML($blankParaXml + $newParaXml)
